$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'241.84"

# Row 3
$ws.Range("D3").Value = "'21.91"

# Row 4
$ws.Range("D4").Value = "'5.386"

# Row 5
$ws.Range("D5").Value = "'0.05713"

# Row 7
$ws.Range("D7").Value = "'6.302"

# Row 8
$ws.Range("D8").Value = "'0.8078"

# Row 9
$ws.Range("D9").Value = "'0.8518"

# Row 10
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = "'0.01089"
$ws.Range("E10").Value = '9OneONEBestin24h'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1438"
$ws.Range("E11").Value = '10WazirXWRX'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.07292"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = "'0.03075"
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = "'0.03142"
$ws.Range("E14").Value = '13BitrueCoinBTR'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.09363"
$ws.Range("E15").Value = '14BitMartTokenBMX'

# Row 16
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = "'3.924"
$ws.Range("E16").Value = '15MCDexMCB'

# Row 17
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = "'0.001587"
$ws.Range("E17").Value = '16BitForexTokenBF'

# Row 18
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = "'0.04811"
$ws.Range("E18").Value = '17CoinExTokenCET'

# Row 19
$ws.Range("D19").Value = "'0.006414"

# Row 20
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").Value = "'0.001001"
$ws.Range("E20").Value = '19BitKanKAN'

# Row 21
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = "'0.004077"
$ws.Range("E21").Value = '20HotbitTokenHTB'

# Row 22
$ws.Range("D22").Value = "'0.0001505"

# Row 23
$ws.Range("D23").Value = "'3.721"

# Row 24
$ws.Range("D24").Value = "'2.172"

# Row 26
$ws.Range("D26").Value = "'0.1304"

# Row 27
$ws.Range("D27").Value = "'0.0004012"

# Row 40
$ws.Range("D40").Value = "'0.03842"

# Row 41
$ws.Range("D41").Value = "'0.006753"

# Row 42
$ws.Range("D42").Value = "'0.1055"

# Row 43
$ws.Range("D43").Value = "'0.002809"

# Row 44
$ws.Range("D44").Value = "'0.007338"

# Row 45
$ws.Range("D45").Value = "'0.00005622"

# Row 47
$ws.Range("D47").Value = "'0.5818"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'

# Row 48
$ws.Range("D48").Value = "'0.1440"

# Row 49
$ws.Range("D49").Value = "'0.00002107"

# Row 50
$ws.Range("D50").Value = "'0.01013"
